$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 ("Report Version" row),
# shifting it down to row 11, then fill the new row 10 with the
# "Gs in Bump Case" / "G_bump" / "Gs" entry.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "Gs in Bump Case"
$ws.Range("B10").Value = "G_bump"
$ws.Range("C10").Value = "Gs"

$ws.Range("C10").Select()
